$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing header cell formatting (bold font, thin border, centered) from H1
# onto the two new header cells so they match the other headers (reuses the same cell style).
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Set the new header labels
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill in the data values for rows 2 through 73 (I and J share the same values)
$values = @(
    8, 7, 7, 9, 9, 9, 9, 9, 9, 9, 9, 9, 9, 9, 9, 9, 11, 9, 10, 9, 9, 9, 9, 9, 9, 8, 9, 9, 9, 9, 9, 10, 9, 9, 9, 9, 7, 11, 8, 8, 8, 9, 8, 8, 8, 9, 9, 9, 9, 9, 9, 7, 9, 9, 9, 9, 9, 9, 9, 9, 9, 9, 9, 9, 9, 6, 6, 5, 5, 3, 3, 2
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $values[$i]
    $ws.Cells.Item($row, 10).Value = $values[$i]
}

